$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Exact "728×9=" "440×7="
Replace-Exact "953×4=" "751×4="
Replace-Exact "809×9=" "883×3="
Replace-Exact "795×2=" "603×8="
Replace-Exact "684×2=" "140×4="
Replace-Exact "956×6=" "811×5="
Replace-Exact "590×7=" "370×7="
Replace-Exact "255×3=" "253×6="
Replace-Exact "405×4=" "223×3="
Replace-Exact "941×7=" "618×7="
Replace-Exact "984×5=" "233×5="
Replace-Exact "168×7=" "737×7="
Replace-Exact "224×4=" "811×7="
Replace-Exact "985×3=" "430×9="
Replace-Exact "365×5=" "642×9="
Replace-Exact "346×3=" "328×5="
Replace-Exact "645×2=" "481×7="
Replace-Exact "648×3=" "400×8="
Replace-Exact "525×4=" "666×7="
Replace-Exact "765×4=" "776×3="
Replace-Exact "401×9=" "997×8="
Replace-Exact "303×6=" "744×2="
Replace-Exact "906×4=" "335×9="
Replace-Exact "218×6=" "328×2="
Replace-Exact "225×3=" "877×5="
